$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("School")

# --- Row 13: new effort entry for 05.11.2015 ---

# A13 holds a date-like label that must stay plain text (like the other
# entries in this "Date" column), even though the cell is formatted with a
# date number format. Entering "05.11.2015" directly would be auto-parsed
# into a real date serial number, so instead we evaluate it via a formula
# that yields the text, then convert that formula to a static value with
# PasteSpecial (values only) - pasted values are not re-parsed as dates.
$ws.Range("A13").Formula2 = '="05.11.2015"'
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range("B13").Value = "Arpan Kar"
$ws.Range("C13").Value = "Coding"
$ws.Range("D13").Value = 4
$ws.Range("J13").Value = "Restructuring UI using control. Training for Sample code."

# Move the active selection to A14, matching the updated sheet view
$ws.Range("A14").Select()

$excel.CalculateFullRebuild()
